# This script adds three new annotation columns ("Parameter [Isotope labeling]",
# "Term Source REF (PRIDE:0000310)", "Term Accession Number (PRIDE:0000310)")
# to the "MassSpec" table, inserted right before the existing "Output [Data]"
# block (which gets shifted three columns to the right), and fills in the
# values for the three data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MassSpec")
$tbl = $ws.ListObjects.Item(1)

# --- Step 1: shift the existing "Output [Data]" / "Data Format " /
# "Data Selector Format " columns (M:O) three places to the right (P:R)
# so that there is room for the three new columns in between.
$ws.Range("M1:O4").Copy()
$ws.Range("P1").PasteSpecial()

# --- Step 2: grow the table to cover the new range A1:R4 (18 columns).
$tbl.Resize($ws.Range("A1:R4"))

# --- Step 3: write the headers for the three new columns into M1:O1.
$ws.Range("M1").Value() = "Parameter [Isotope labeling]"
$ws.Range("N1").Value() = "Term Source REF (PRIDE:0000310)"
$ws.Range("O1").Value() = "Term Accession Number (PRIDE:0000310)"

# Re-assert the (already shifted) headers in P1:R1 so the table's column
# metadata picks up their names instead of generic placeholders.
$ws.Range("P1").Value() = $ws.Range("P1").Value()
$ws.Range("Q1").Value() = $ws.Range("Q1").Value()
$ws.Range("R1").Value() = $ws.Range("R1").Value()

# --- Step 4: populate the new columns for each of the three data rows.
$ws.Range("M2").Value() = "metabolic labelling: heavy N (mainly 15N)"
$ws.Range("N2").Value() = "MS"
$ws.Range("O2").Value() = "https://www.ebi.ac.uk/ols4/ontologies/ms/classes/http%253A%252F%252Fpurl.obolibrary.org%252Fobo%252FMS_1002068"

$ws.Range("M3").Value() = "metabolic labelling: heavy N (mainly 15N)"
$ws.Range("N3").Value() = "MS"
$ws.Range("O3").Value() = "https://www.ebi.ac.uk/ols4/ontologies/ms/classes/http%253A%252F%252Fpurl.obolibrary.org%252Fobo%252FMS_1002068"

$ws.Range("M4").Value() = "metabolic labelling: heavy N (mainly 15N)"
$ws.Range("N4").Value() = "MS"
$ws.Range("O4").Value() = "https://www.ebi.ac.uk/ols4/ontologies/ms/classes/http%253A%252F%252Fpurl.obolibrary.org%252Fobo%252FMS_1002068"

Write-Host "Finished adding Isotope labeling columns to MassSpec table."
